# Debugging the Cloud run
#
# GEPEP_calibration.xlsx - "Low" scenario sheet:
#   - Column I (PI) rows 2-31: 0.03 -> 0.035
#   - Column K (Wwt scalar) rows 2-31: 0 -> 0.00000009, displayed in
#     scientific notation (0.00E+00)
#   - The "Low" sheet tab becomes the active/selected sheet (it was
#     "BestBet" before), with K2:K31 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Low")

# --- Column I: PI, 0.03 -> 0.035 for rows 2-31 --------------------------
$ws.Range("I2:I31").Value = 0.035

# --- Column K: Wwt scalar, 0 -> 9E-08, shown in scientific notation ----
$ws.Range("K2:K31").Value = 0.00000009
$ws.Range("K2:K31").NumberFormat = "0.00E+00"

# --- View state: make "Low" the active sheet/tab, with K2:K31 selected -
$ws.Activate()
$ws.Range("K2:K31").Select() | Out-Null
